$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 254:262 column D: stray inline-string "0.00" -> numeric 0 ---
for ($r = 254; $r -le 262; $r++) {
  $ws.Cells.Item($r, 4).Value = 0
}

# --- Append new traffic rows 263:300 ---
$newRows = @(
  @{Row=263; A='2025-06-03 20:41:24'; B=2; C='Car'; D='31.51'},
  @{Row=264; A='2025-06-03 20:41:25'; B=4; C='Car'; D='45.77'},
  @{Row=265; A='2025-06-03 20:41:25'; B=5; C='Car'; D='41.46'},
  @{Row=266; A='2025-06-03 20:41:31'; B=9; C='Car'; D='41.01'},
  @{Row=267; A='2025-06-03 20:41:33'; B=10; C='Truck'; D='51.41'},
  @{Row=268; A='2025-06-03 20:41:34'; B=13; C='Car'; D='41.89'},
  @{Row=269; A='2025-06-03 20:41:36'; B=14; C='Car'; D='42.43'},
  @{Row=270; A='2025-06-03 20:44:08'; B=25; C='Car'; D='0.00'},
  @{Row=271; A='2025-06-03 20:44:09'; B=27; C='Car'; D='0.00'},
  @{Row=272; A='2025-06-03 20:44:10'; B=28; C='Car'; D='0.00'},
  @{Row=273; A='2025-06-03 20:44:11'; B=31; C='Car'; D='0.00'},
  @{Row=274; A='2025-06-03 20:44:14'; B=34; C='Bus'; D='0.00'},
  @{Row=275; A='2025-06-03 20:44:15'; B=32; C='Car'; D='0.00'},
  @{Row=276; A='2025-06-03 20:44:16'; B=35; C='Car'; D='0.00'},
  @{Row=277; A='2025-06-03 20:44:18'; B=33; C='Truck'; D='0.00'},
  @{Row=278; A='2025-06-03 20:44:18'; B=36; C='Car'; D='0.00'},
  @{Row=279; A='2025-06-03 20:44:21'; B=37; C='Car'; D='0.00'},
  @{Row=280; A='2025-06-03 20:44:21'; B=39; C='Car'; D='0.00'},
  @{Row=281; A='2025-06-03 20:44:23'; B=40; C='Car'; D='0.00'},
  @{Row=282; A='2025-06-03 20:44:24'; B=41; C='Car'; D='0.00'},
  @{Row=283; A='2025-06-03 20:44:28'; B=49; C='Truck'; D='38.01'},
  @{Row=284; A='2025-06-03 20:44:29'; B=50; C='Car'; D='43.20'},
  @{Row=285; A='2025-06-03 20:44:30'; B=51; C='Truck'; D='48.93'},
  @{Row=286; A='2025-06-03 20:44:31'; B=46; C='Car'; D='90.28'},
  @{Row=287; A='2025-06-03 20:44:33'; B=52; C='Bus'; D='15.26'},
  @{Row=288; A='2025-06-03 20:44:35'; B=53; C='Car'; D='44.99'},
  @{Row=289; A='2025-06-03 20:44:36'; B=54; C='Car'; D='42.89'},
  @{Row=290; A='2025-06-03 20:44:38'; B=55; C='Car'; D='40.32'},
  @{Row=291; A='2025-06-03 20:44:38'; B=56; C='Car'; D='44.76'},
  @{Row=292; A='2025-06-03 20:44:41'; B=57; C='Car'; D='35.97'},
  @{Row=293; A='2025-06-03 20:44:41'; B=58; C='Car'; D='41.35'},
  @{Row=294; A='2025-06-03 20:44:43'; B=61; C='Car'; D='51.75'},
  @{Row=295; A='2025-06-03 20:44:51'; B=64; C='Truck'; D='123.76'},
  @{Row=296; A='2025-06-03 20:44:54'; B=66; C='Car'; D='87.57'},
  @{Row=297; A='2025-06-03 20:44:55'; B=68; C='Car'; D='74.44'},
  @{Row=298; A='2025-06-03 20:44:57'; B=67; C='Car'; D='112.09'},
  @{Row=299; A='2025-06-03 20:44:58'; B=69; C='Car'; D='111.01'},
  @{Row=300; A='2025-06-03 20:45:02'; B=70; C='Truck'; D='19.27'}
)

foreach ($nr in $newRows) {
  $ws.Cells.Item($nr.Row, 1).Value = $nr.A
  $ws.Cells.Item($nr.Row, 2).Value = $nr.B
  $ws.Cells.Item($nr.Row, 3).Value = $nr.C
  $dCell = $ws.Cells.Item($nr.Row, 4)
  $dCell.NumberFormat = "@"
  $dCell.Value = $nr.D
  $dCell.ClearFormats()
}
